# Insert a new weekly price record as row 4, pushing the existing rows
# (4..112) down to (5..113), matching the authors "Fruta / hortaliza,
# semanal" weekly-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:112 down to 5:113 by inserting a fresh row at position 4.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the new daily/weekly price entry.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45160
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112001
$ws.Range("G4").Value = "Berenjena"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 9000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 9000
$ws.Range("N4").Value = "$/caja 60 unidades"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 150
$ws.Range("Q4").Value = 60
$ws.Range("R4").Value = "Hortaliza"
